# OrderReceipts_Expenses_Konto_Classification.xlsx
# Entered 18 new rows (1340-1357) on the "Konto" sheet for week ending
# 2021-07-11 (total 1356 data rows + header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1340 ---
$ws.Cells.Item(1340, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1340, 1).Value = 44384
$ws.Cells.Item(1340, 2).Value = 'Reko253'
$ws.Cells.Item(1340, 3).Value = 3011
$ws.Cells.Item(1340, 4).Value = 'Reko Swish +46733496329'
$ws.Cells.Item(1340, 6).Value = 141.07

# --- Row 1341 ---
$ws.Cells.Item(1341, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1341, 1).Value = 44384
$ws.Cells.Item(1341, 2).Value = 'Reko253'
$ws.Cells.Item(1341, 3).Value = 2611
$ws.Cells.Item(1341, 4).Value = 'Reko Swish +46733496329'
$ws.Cells.Item(1341, 6).Value = 16.93

# --- Row 1342 ---
$ws.Cells.Item(1342, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1342, 1).Value = 44384
$ws.Cells.Item(1342, 2).Value = 'Reko253'
$ws.Cells.Item(1342, 3).Value = 1930
$ws.Cells.Item(1342, 4).Value = 'Reko Swish +46733496329'
$ws.Cells.Item(1342, 5).Value = 158

# --- Row 1343 ---
$ws.Cells.Item(1343, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1343, 1).Value = 44384
$ws.Cells.Item(1343, 2).Value = 'Reko254'
$ws.Cells.Item(1343, 3).Value = 3011
$ws.Cells.Item(1343, 4).Value = 'Reko Swish +46702597315'
$ws.Cells.Item(1343, 6).Value = 282.14

# --- Row 1344 ---
$ws.Cells.Item(1344, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1344, 1).Value = 44384
$ws.Cells.Item(1344, 2).Value = 'Reko254'
$ws.Cells.Item(1344, 3).Value = 2611
$ws.Cells.Item(1344, 4).Value = 'Reko Swish +46702597315'
$ws.Cells.Item(1344, 6).Value = 33.86

# --- Row 1345 ---
$ws.Cells.Item(1345, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1345, 1).Value = 44384
$ws.Cells.Item(1345, 2).Value = 'Reko254'
$ws.Cells.Item(1345, 3).Value = 1930
$ws.Cells.Item(1345, 4).Value = 'Reko Swish +46702597315'
$ws.Cells.Item(1345, 5).Value = 316

# --- Row 1346 ---
$ws.Cells.Item(1346, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1346, 1).Value = 44385
$ws.Cells.Item(1346, 2).Value = 'Reko255'
$ws.Cells.Item(1346, 3).Value = 3011
$ws.Cells.Item(1346, 4).Value = 'Reko Swish +46707255040'
$ws.Cells.Item(1346, 6).Value = 211.61

# --- Row 1347 ---
$ws.Cells.Item(1347, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1347, 1).Value = 44385
$ws.Cells.Item(1347, 2).Value = 'Reko255'
$ws.Cells.Item(1347, 3).Value = 2611
$ws.Cells.Item(1347, 4).Value = 'Reko Swish +46707255040'
$ws.Cells.Item(1347, 6).Value = 25.39

# --- Row 1348 ---
$ws.Cells.Item(1348, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1348, 1).Value = 44385
$ws.Cells.Item(1348, 2).Value = 'Reko255'
$ws.Cells.Item(1348, 3).Value = 1930
$ws.Cells.Item(1348, 4).Value = 'Reko Swish +46707255040'
$ws.Cells.Item(1348, 5).Value = 237

# --- Row 1349 ---
$ws.Cells.Item(1349, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1349, 1).Value = 44385
$ws.Cells.Item(1349, 2).Value = 'Reko256'
$ws.Cells.Item(1349, 3).Value = 3011
$ws.Cells.Item(1349, 4).Value = 'Reko Swish +46723030040'
$ws.Cells.Item(1349, 6).Value = 115.18

# --- Row 1350 ---
$ws.Cells.Item(1350, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1350, 1).Value = 44385
$ws.Cells.Item(1350, 2).Value = 'Reko256'
$ws.Cells.Item(1350, 3).Value = 2611
$ws.Cells.Item(1350, 4).Value = 'Reko Swish +46723030040'
$ws.Cells.Item(1350, 6).Value = 13.82

# --- Row 1351 ---
$ws.Cells.Item(1351, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1351, 1).Value = 44385
$ws.Cells.Item(1351, 2).Value = 'Reko256'
$ws.Cells.Item(1351, 3).Value = 1930
$ws.Cells.Item(1351, 4).Value = 'Reko Swish +46723030040'
$ws.Cells.Item(1351, 5).Value = 129

# --- Row 1352 ---
# Receipt number "5081937" looks numeric but must stay a text value (as in the
# source data) - enter it with a leading apostrophe (forces text, same as
# typing it into Excel) and then reset the resulting "quote prefix" cell
# style back to Normal so no visible/number formatting change sticks.
$ws.Cells.Item(1352, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1352, 1).Value = 44385
$ws.Cells.Item(1352, 2).Value = "'5081937"
$ws.Cells.Item(1352, 2).Style = "Normal"
$ws.Cells.Item(1352, 3).Value = 3011
$ws.Cells.Item(1352, 4).Value = 'Order 5081937 Swish +46703564388'
$ws.Cells.Item(1352, 6).Value = 1062.5

# --- Row 1353 ---
$ws.Cells.Item(1353, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1353, 1).Value = 44385
$ws.Cells.Item(1353, 2).Value = "'5081937"
$ws.Cells.Item(1353, 2).Style = "Normal"
$ws.Cells.Item(1353, 3).Value = 2611
$ws.Cells.Item(1353, 4).Value = 'Order 5081937 Swish +46703564388'
$ws.Cells.Item(1353, 6).Value = 127.5

# --- Row 1354 ---
$ws.Cells.Item(1354, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1354, 1).Value = 44385
$ws.Cells.Item(1354, 2).Value = "'5081937"
$ws.Cells.Item(1354, 2).Style = "Normal"
$ws.Cells.Item(1354, 3).Value = 1930
$ws.Cells.Item(1354, 4).Value = 'Order 5081937 Swish +46703564388'
$ws.Cells.Item(1354, 5).Value = 1190

# --- Row 1355 ---
$ws.Cells.Item(1355, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1355, 1).Value = 44388
$ws.Cells.Item(1355, 3).Value = 5010
$ws.Cells.Item(1355, 4).Value = 'July hyra'
$ws.Cells.Item(1355, 5).Value = 4166

# --- Row 1356 ---
$ws.Cells.Item(1356, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1356, 1).Value = 44388
$ws.Cells.Item(1356, 4).Value = 'July hyra'
$ws.Cells.Item(1356, 5).Value = 0

# --- Row 1357 ---
$ws.Cells.Item(1357, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1357, 1).Value = 44388
$ws.Cells.Item(1357, 3).Value = 1930
$ws.Cells.Item(1357, 4).Value = 'July hyra'
$ws.Cells.Item(1357, 6).Value = 4166
